{"js": "// Update the date line (first paragraph, above the table) and the 5x5\n// grid of three-digit \u00f7 one-digit division problems in the table below\n// it. Cell text is replaced via the paragraph's own Range so existing\n// run/paragraph formatting (font, size, alignment) is preserved.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].getRange().insertText(\"2025-01-31 Friday\", Word.InsertLocation.replace);\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row/column indices are 0-based; the worked rows in the table are\n// 0, 4, 8, 12, 16 (every row in between is a blank spacer row).\nconst newValues = [\n  [0, [\"455\u00f77=65, 0\", \"380\u00f79=42, 2\", \"634\u00f78=79, 2\", \"246\u00f79=27, 3\", \"290\u00f79=32, 2\"]],\n  [4, [\"264\u00f74=66, 0\", \"643\u00f77=91, 6\", \"840\u00f74=210, 0\", \"195\u00f73=65, 0\", \"313\u00f73=104, 1\"]],\n  [8, [\"601\u00f75=120, 1\", \"370\u00f76=61, 4\", \"858\u00f75=171, 3\", \"407\u00f78=50, 7\", \"316\u00f79=35, 1\"]],\n  [12, [\"473\u00f76=78, 5\", \"544\u00f72=272, 0\", \"568\u00f76=94, 4\", \"930\u00f72=465, 0\", \"467\u00f76=77, 5\"]],\n  [16, [\"440\u00f75=88, 0\", \"607\u00f72=303, 1\", \"636\u00f74=159, 0\", \"417\u00f78=52, 1\", \"631\u00f77=90, 1\"]],\n];\n\nfor (const [rowIndex, rowValues] of newValues) {\n  for (let colIndex = 0; colIndex < rowValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n\n    cellParagraphs.items[0].getRange().insertText(rowValues[colIndex], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph, above the table)\n$d.Paragraphs.Item(1).Range.Text = \"2025-01-31 Friday\"\n\n# Update the 5x5 grid of division problems that live in rows 1, 5, 9, 13, 17\n$t = $d.Tables.Item(1)\n\n$values = @(\n  @(\"455\u00f77=65, 0\", \"380\u00f79=42, 2\", \"634\u00f78=79, 2\", \"246\u00f79=27, 3\", \"290\u00f79=32, 2\"),\n  @(\"264\u00f74=66, 0\", \"643\u00f77=91, 6\", \"840\u00f74=210, 0\", \"195\u00f73=65, 0\", \"313\u00f73=104, 1\"),\n  @(\"601\u00f75=120, 1\", \"370\u00f76=61, 4\", \"858\u00f75=171, 3\", \"407\u00f78=50, 7\", \"316\u00f79=35, 1\"),\n  @(\"473\u00f76=78, 5\", \"544\u00f72=272, 0\", \"568\u00f76=94, 4\", \"930\u00f72=465, 0\", \"467\u00f76=77, 5\"),\n  @(\"440\u00f75=88, 0\", \"607\u00f72=303, 1\", \"636\u00f74=159, 0\", \"417\u00f78=52, 1\", \"631\u00f77=90, 1\")\n)\n\n$rows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $rows.Length; $i++) {\n  $r = $rows[$i]\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $values[$i][$c - 1]\n  }\n}\n"}
